$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 12
$ws.Range("I31").Value = 12
$ws.Range("K31").Value = 36
$ws.Range("M31").Value = 194

$ws.Range("H58").Value = 1879.875
$ws.Range("I58").Value = 805
$ws.Range("K58").Value = 2415
$ws.Range("M58").Value = -2265

$ws.Range("H70").Value = 1067695.6
$ws.Range("I70").Value = 1841452.9
$ws.Range("K70").Value = 5524358.699999999
$ws.Range("M70").Value = -5524088.699999999

$ws.Range("H73").Value = 1067695.6
$ws.Range("I73").Value = 1841452.9
$ws.Range("K73").Value = 5524358.699999999
$ws.Range("M73").Value = -5523422.699999999

$ws.Range("H92").Value = 583.3333
$ws.Range("I92").Value = 650
$ws.Range("K92").Value = 650
$ws.Range("M92").Value = 598

$ws.Range("H107").Value = 1925.6666
$ws.Range("I107").Value = 925.5
$ws.Range("K107").Value = 925.5
$ws.Range("M107").Value = 994.5

$ws.Range("H111").Value = 753
$ws.Range("I111").Value = 753
$ws.Range("K111").Value = 2259
$ws.Range("M111").Value = 808

$ws.Range("H125").Value = 7026.5454
$ws.Range("I125").Value = 5286.625
$ws.Range("K125").Value = 47579.625
$ws.Range("M125").Value = -45119.625

$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140

$ws.Range("H135").Value = 1997.5
$ws.Range("I135").Value = 1496.75
$ws.Range("K135").Value = 13470.75
$ws.Range("M135").Value = -10935.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8830.939
$ws.Range("I32").Value = 7544.4062
$ws.Range("K32").Value = 7544.4062
$ws.Range("M32").Value = -7257.4062

$ws.Range("H61").Value = 1800
$ws.Range("J61").Value = 1800
$ws.Range("L61").Value = 1800
$ws.Range("N61").Value = -2224

$ws.Range("H74").Value = 2426.875
$ws.Range("I74").Value = 2750.25
$ws.Range("J74").Value = 2103.5
$ws.Range("K74").Value = 2750.25
$ws.Range("L74").Value = 2103.5
$ws.Range("M74").Value = -1876.25
$ws.Range("N74").Value = -3851.5

$ws.Range("H77").Value = 2426.875
$ws.Range("I77").Value = 2750.25
$ws.Range("J77").Value = 2103.5
$ws.Range("K77").Value = 13751.25
$ws.Range("L77").Value = 10517.5
$ws.Range("M77").Value = -9383.25
$ws.Range("N77").Value = -19253.5

$ws.Range("H122").Value = 3416.1428
$ws.Range("I122").Value = 3416.1428
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10248.4284
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7798.428400000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4662.846
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1800
$ws.Range("J136").Value = 1800
$ws.Range("L136").Value = 5400
$ws.Range("N136").Value = -10500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1399.92
$ws.Range("I20").Value = 1181.95
$ws.Range("J20").Value = 2271.8
$ws.Range("K20").Value = 1181.95
$ws.Range("L20").Value = 2271.8
$ws.Range("M20").Value = -934.95
$ws.Range("N20").Value = -2765.8

$ws.Range("H86").Value = 2370.0293
$ws.Range("I86").Value = 2273.963
$ws.Range("K86").Value = 2273.963
$ws.Range("M86").Value = -1150.963

$ws.Range("H89").Value = 2370.0293
$ws.Range("I89").Value = 2273.963
$ws.Range("K89").Value = 11369.815
$ws.Range("M89").Value = -5753.815000000001

$ws.Range("H94").Value = 1352.48
$ws.Range("I94").Value = 1390.3334
$ws.Range("J94").Value = 1153.75
$ws.Range("K94").Value = 1390.3334
$ws.Range("L94").Value = 1153.75
$ws.Range("M94").Value = -939.3334
$ws.Range("N94").Value = -2055.75

$ws.Range("H99").Value = 3799.4
$ws.Range("I99").Value = 3799.4
$ws.Range("K99").Value = 3799.4
$ws.Range("M99").Value = -2301.4

$ws.Range("H134").Value = 5422.36
$ws.Range("I134").Value = 5481.625
$ws.Range("K134").Value = 16444.875
$ws.Range("M134").Value = -13909.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 99
$ws.Range("I4").Value = 99
$ws.Range("K4").Value = 99
$ws.Range("M4").Value = 13

$ws.Range("H11").Value = 686.3333
$ws.Range("J11").Value = 537.5
$ws.Range("L11").Value = 537.5
$ws.Range("N11").Value = -817.5

$ws.Range("H103").Value = 7075
$ws.Range("I103").Value = 7075
$ws.Range("K103").Value = 7075
$ws.Range("M103").Value = -5903

$ws.Range("H122").Value = 4210
$ws.Range("I122").Value = 5041.4287
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 15124.2861
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -12674.2861
$ws.Range("N122").Value = -8800

$ws.Range("H134").Value = 2349.9524
$ws.Range("I134").Value = 1620.9412
$ws.Range("K134").Value = 4862.8236
$ws.Range("M134").Value = -2327.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 464.5
$ws.Range("I23").Value = 546.6667
$ws.Range("J23").Value = 415.2
$ws.Range("K23").Value = 1640.0001
$ws.Range("L23").Value = 1245.6
$ws.Range("M23").Value = -1405.0001
$ws.Range("N23").Value = -1715.6

$ws.Range("H36").Value = 166.33333
$ws.Range("I36").Value = 166.33333
$ws.Range("K36").Value = 498.99999
$ws.Range("M36").Value = -329.99999

$ws.Range("H40").Value = 180.82353
$ws.Range("I40").Value = 144.15384
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 576.61536
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -507.61536
$ws.Range("N40").Value = -1338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10302

$ws.Range("H97").Value = 804.1875
$ws.Range("I97").Value = 849.4666999999999
$ws.Range("J97").Value = 125
$ws.Range("K97").Value = 849.4666999999999
$ws.Range("L97").Value = 125
$ws.Range("M97").Value = -353.4666999999999
$ws.Range("N97").Value = -1117

$ws.Range("H102").Value = 1536.6154
$ws.Range("I102").Value = 1536.6154
$ws.Range("K102").Value = 1536.6154
$ws.Range("M102").Value = 85.38460000000009

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 5999.3335
$ws.Range("J13").Value = 5999.3335
$ws.Range("L13").Value = 5999.3335
$ws.Range("N13").Value = -6279.3335

$ws.Range("H16").Value = 1933.2
$ws.Range("I16").Value = 1933.2
$ws.Range("K16").Value = 1933.2
$ws.Range("M16").Value = -1763.2

$ws.Range("H26").Value = 1649.5
$ws.Range("J26").Value = 1649.5
$ws.Range("L26").Value = 1649.5
$ws.Range("N26").Value = -2239.5

$ws.Range("H61").Value = 6215.6665
$ws.Range("I61").Value = 5159.3
$ws.Range("K61").Value = 5159.3
$ws.Range("M61").Value = -4957.3

$ws.Range("H106").Value = 12996.667
$ws.Range("J106").Value = 12996.667
$ws.Range("L106").Value = 12996.667
$ws.Range("N106").Value = -15520.667

$ws.Range("H113").Value = 6215.6665
$ws.Range("I113").Value = 5159.3
$ws.Range("K113").Value = 5159.3
$ws.Range("M113").Value = -2989.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 731.5
$ws.Range("I100").Value = 692.25
$ws.Range("K100").Value = 1384.5
$ws.Range("M100").Value = -843.5

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12530
$ws.Range("N126").ClearContents()

$ws.Range("H135").Value = 48411.555
$ws.Range("J135").Value = 52213
$ws.Range("L135").Value = 52213
$ws.Range("N135").Value = -62353
